# Apply Swahili -> English translation edits via Find + direct Range.Text assignment
# (avoids Find.Execute's ReplaceWith smart-quote autocorrect substitution)
$d = $word.ActiveDocument

$pairs = @(
    ,@('Umbizo limesahihishwa sio wakati', 'Format has been corrected not the timing')
    ,@('Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino', 'I added 25 seconds to each timing to correct for the intro song -john argentino')
    ,@('Tatizo la uwanja wa ndege - manukuu:', 'The airport problem - subtitles:')
    ,@('Utawala wa tatu', 'The administrations of three')
    ,@('miji jirani: A, B na C waliamua', 'neighboring cities: A, B and C decided')
    ,@('kujenga uwanja wa ndege unaogawanya gharama za', 'to build an airport dividing the costs of')
    ,@('utekelezaji. Hali juu ya', 'implementation. The condition on the')
    ,@('uchaguzi wa mahali pa kufaa zaidi ni', 'choice of the most suitable place is')
    ,@('kwamba jumla ya umbali kutoka kwa kila mmoja', 'that the sum of the distances from each')
    ,@('mji kwa uwanja wa ndege ni ndogo kama', 'city to the airport is as small as')
    ,@('inawezekana. Timu ya wataalam wanaohusika', 'possible. The team of experts in charge')
    ,@('ya kazi imeunda mfano wa kupata', 'of the work has created a model to get')
    ,@('wazo la awali la mahali pa kuweka', 'a preliminary idea of where to place the')
    ,@('muundo. Ovyo wao wapo', 'structure. At their disposal there are')
    ,@('konokono wengine pete kubwa ya chuma na ndefu', 'some snails a big metal ring and a long')
    ,@('kamba.', 'string.')
    ,@('Eleza jinsi timu inaweza kusimamia matumizi', 'Explain how the team can manage to use')
    ,@('nyenzo za kusema takriban', 'the materials to tell approximately the')
    ,@('eneo bora la uwanja wa ndege. Fikiria', 'ideal location of the airport. Imagine')
    ,@('kwamba miji imewekwa kwenye', 'that the cities are placed at the')
    ,@('vipeo vya pembetatu ambayo ni', 'vertices of a triangle which is')
    ,@('kwa hakika imetolewa tena kwa kiwango kama', 'obviously reproduced in scale as')
    ,@('inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana', 'shown in figure. This is one possible')
    ,@('kuweka kamba huanza kutoka msumari mmoja,', 'setting the rope starts from one nail,')
    ,@('huenda ndani ya pete, huzunguka', 'goes inside the ring, goes around the')
    ,@('msumari mwingine, msumari wa tatu, ndani ya', 'other nail, the third nail, inside the')
    ,@('pete tena na sasa unaweza kuvuta tu', 'ring again and now you can just pull the')
    ,@('kamba ili kupata uhakika huo', 'rope in order to find the point that')
    ,@('unatafuta. Ili kufikia', 'you''re looking for. In order to reach the')
    ,@('uhakika, tunapaswa kusonga kamba kidogo', 'point, we have to move the rope a bit')
    ,@('kwa sababu kuna ', 'because there is some ')
    ,@('upinzani', 'resistance')
    ,@(' uliosababishwa', ' caused')
    ,@('kwa nyenzo ambazo tunatumia lakini', 'by the materials that we are using but')
    ,@('baada ya muda utafikia nafasi kutoka', 'after a while you''ll reach a position from')
    ,@('ambayo pete haisogei tena,', 'which the ring doesn''t move anymore,')
    ,@('ambayo ni zaidi au chini ya hii. Na kama', 'which is more or less this one. And as')
    ,@('kati ya pete na misumari ni', 'between the ring and the nails are')
    ,@('kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja', 'placed more or less 120 degrees from one')
    ,@('nyingine ambayo ni 1/3 ya mduara,', 'another which is 1/3 of a circumference,')
    ,@('na hiyo ndiyo hatua tunayoiangalia', 'and that''s the point that we''re looking')
    ,@('kwa: umbali wa chini kati ya', 'for: the minimum distance between the')
    ,@('misumari na uwanja wa ndege unapojumlisha', 'nails and the airport when you sum it')
    ,@('pamoja', 'ogether')
    ,@('[Muziki]', '[Music]')
)

foreach ($pair in $pairs) {
    $findText = $pair[0]
    $replText = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $range.Text = $replText
    } else {
        Write-Output "NOT FOUND: $findText"
    }
}

Write-Output "Done"
